$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 170, pushing the existing rows 170-296 down to 172-298.
$ws.Rows("170:171").Insert()

# Row 170 (new): "Primera" quality record for 2022-01-01 (serial 44566)
$ws.Cells.Item(170,1).Value2  = 8
$ws.Cells.Item(170,2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(170,3).Value2  = "Coquimbo"
$ws.Cells.Item(170,4).Value2  = 44566
$ws.Cells.Item(170,5).Value2  = 4
$ws.Cells.Item(170,6).Value2  = 100112017
$ws.Cells.Item(170,7).Value2  = "Apio"
$ws.Cells.Item(170,8).Value2  = "Americana (o)"
$ws.Cells.Item(170,9).Value2  = "Primera"
$ws.Cells.Item(170,10).Value2 = 2500
$ws.Cells.Item(170,11).Value2 = 8000
$ws.Cells.Item(170,12).Value2 = 9000
$ws.Cells.Item(170,13).Value2 = 8500
$ws.Cells.Item(170,14).Value2 = "`$/docena de matas"
$ws.Cells.Item(170,15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(170,16).Value2 = 1417
$ws.Cells.Item(170,17).Value2 = 6
$ws.Cells.Item(170,18).Value2 = "Hortaliza"

# Row 171 (new): "Segunda" quality record for the same date
$ws.Cells.Item(171,1).Value2  = 8
$ws.Cells.Item(171,2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(171,3).Value2  = "Coquimbo"
$ws.Cells.Item(171,4).Value2  = 44566
$ws.Cells.Item(171,5).Value2  = 4
$ws.Cells.Item(171,6).Value2  = 100112017
$ws.Cells.Item(171,7).Value2  = "Apio"
$ws.Cells.Item(171,8).Value2  = "Americana (o)"
$ws.Cells.Item(171,9).Value2  = "Segunda"
$ws.Cells.Item(171,10).Value2 = 1480
$ws.Cells.Item(171,11).Value2 = 6000
$ws.Cells.Item(171,12).Value2 = 7000
$ws.Cells.Item(171,13).Value2 = 6500
$ws.Cells.Item(171,14).Value2 = "`$/docena de matas"
$ws.Cells.Item(171,15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(171,16).Value2 = 1083
$ws.Cells.Item(171,17).Value2 = 6
$ws.Cells.Item(171,18).Value2 = "Hortaliza"
